$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-column (D) cells being updated below hold text values such as
# "1.00" / "62.763.82" (not real numbers). Force text format on just those
# cells first so Excel does not silently coerce them to numbers/dates.
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D10", "D13", "D15", "D17", "D18", "D19", "D21", "D22", "D24", "D25", "D26", "D28", "D29", "D31", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D47", "D49", "D50", "D51")
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '62.763.82'
$ws.Range("E2").Value = '  -2.53%  '
$ws.Range("D3").Value = '3.392.13'
$ws.Range("E3").Value = '  -3.61%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '574.23'
$ws.Range("E5").Value = '  -2.98%  '
$ws.Range("D6").Value = '125.82'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.393.83'
$ws.Range("E8").Value = '  -3.53%  '
$ws.Range("D9").Value = '0.474'
$ws.Range("E9").Value = '  -2.94%  '
$ws.Range("D10").Value = '7.35'
$ws.Range("E10").Value = '  -3.56%  '
$ws.Range("E11").Value = '  -3.00%  '
$ws.Range("E12").Value = '  -2.68%  '
$ws.Range("D13").Value = '3.967.51'
$ws.Range("E13").Value = '  -3.69%  '
$ws.Range("E14").Value = '  -0.75%  '
$ws.Range("D15").Value = '3.390.11'
$ws.Range("E15").Value = '  -3.65%  '
$ws.Range("E16").Value = '  -4.69%  '
$ws.Range("D17").Value = '62.727.11'
$ws.Range("E17").Value = '  -2.58%  '
$ws.Range("D18").Value = '24.78'
$ws.Range("E18").Value = '  -4.39%  '
$ws.Range("D19").Value = '9.51'
$ws.Range("E19").Value = '  -4.79%  '
$ws.Range("E20").Value = '  -1.52%  '
$ws.Range("D21").Value = '13.17'
$ws.Range("E21").Value = '  -3.08%  '
$ws.Range("D22").Value = '376.60'
$ws.Range("E22").Value = '  -4.53%  '
$ws.Range("E23").Value = '  -3.37%  '
$ws.Range("D24").Value = '3.527.67'
$ws.Range("E24").Value = '  -3.62%  '
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  +0.11%  '
$ws.Range("D26").Value = '72.21'
$ws.Range("E26").Value = '  -3.36%  '
$ws.Range("E27").Value = '  -8.32%  '
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").Value = '6.98'
$ws.Range("E29").Value = '  -5.78%  '
$ws.Range("E30").Value = '  -4.74%  '
$ws.Range("D31").Value = '7.85'
$ws.Range("E31").Value = '  -5.60%  '
$ws.Range("E32").Value = '  -4.83%  '
$ws.Range("E33").Value = '  -4.82%  '
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").Value = '3.420.12'
$ws.Range("E35").Value = '  -3.64%  '
$ws.Range("D36").Value = '22.73'
$ws.Range("E36").Value = '  -3.07%  '
$ws.Range("D37").Value = '5.28'
$ws.Range("E37").Value = '  -1.62%  '
$ws.Range("D38").Value = '6.73'
$ws.Range("E38").Value = '  -3.51%  '
$ws.Range("D39").Value = '164.54'
$ws.Range("E39").Value = '  -1.58%  '
$ws.Range("E40").Value = '  -4.86%  '
$ws.Range("D41").Value = '0.0758'
$ws.Range("E41").Value = '  -4.29%  '
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("D43").Value = '0.774'
$ws.Range("E43").Value = '  -4.71%  '
$ws.Range("D44").Value = '41.50'
$ws.Range("E44").Value = '  -2.00%  '
$ws.Range("D45").Value = '4.27'
$ws.Range("E45").Value = '  -4.08%  '
$ws.Range("E46").Value = '  -5.33%  '
$ws.Range("D47").Value = '22.92'
$ws.Range("E47").Value = '  -10.20%  '
$ws.Range("E48").Value = '  -8.37%  '
$ws.Range("D49").Value = '6.64'
$ws.Range("E49").Value = '  -2.47%  '
$ws.Range("D50").Value = '2.239.12'
$ws.Range("E50").Value = '  -6.93%  '
$ws.Range("D51").Value = '0.853'
$ws.Range("E51").Value = '  -5.18%  '
